$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 16771.334
$ws.Range("I9").Value = 20105.6
$ws.Range("K9").Value = 20105.6
$ws.Range("M9").Value = -19936.6
$ws.Range("H48").Value = 2474.5
$ws.Range("I48").Value = 2450
$ws.Range("K48").Value = 7350
$ws.Range("M48").Value = -7058
$ws.Range("H56").Value = 2474.5
$ws.Range("I56").Value = 2450
$ws.Range("K56").Value = 7350
$ws.Range("M56").Value = -6816
$ws.Range("H82").Value = 4638.222
$ws.Range("I82").Value = 2593.125
$ws.Range("J82").Value = 20999
$ws.Range("K82").Value = 7779.375
$ws.Range("L82").Value = 62997
$ws.Range("M82").Value = -7373.375
$ws.Range("N82").Value = -63809
$ws.Range("H85").Value = 4638.222
$ws.Range("I85").Value = 2593.125
$ws.Range("J85").Value = 20999
$ws.Range("K85").Value = 7779.375
$ws.Range("L85").Value = 62997
$ws.Range("M85").Value = -6375.375
$ws.Range("N85").Value = -65805
$ws.Range("H99").Value = 3209.5715
$ws.Range("I99").Value = 257.4
$ws.Range("J99").Value = 10590
$ws.Range("K99").Value = 772.1999999999999
$ws.Range("L99").Value = 31770
$ws.Range("M99").Value = 725.8000000000001
$ws.Range("N99").Value = -34766
$ws.Range("H100").Value = 2100.0527
$ws.Range("I100").Value = 1820.4667
$ws.Range("K100").Value = 1820.4667
$ws.Range("M100").Value = -1279.4667
$ws.Range("H137").Value = 1308.1578
$ws.Range("I137").Value = 1201.3334
$ws.Range("J137").Value = 1708.75
$ws.Range("K137").Value = 3604.0002
$ws.Range("L137").Value = 5126.25
$ws.Range("M137").Value = -1054.0002
$ws.Range("N137").Value = -10226.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3767.4
$ws.Range("I2").Value = 5594
$ws.Range("J2").Value = 2549.6667
$ws.Range("K2").Value = 5594
$ws.Range("L2").Value = 2549.6667
$ws.Range("M2").Value = -5481
$ws.Range("N2").Value = -2775.6667
$ws.Range("H4").Value = 399.6
$ws.Range("I4").Value = 399.5
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 399.5
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -283.5
$ws.Range("N4").Value = -632
$ws.Range("H5").Value = 3244.25
$ws.Range("I5").Value = 1406
$ws.Range("K5").Value = 1406
$ws.Range("M5").Value = -1294
$ws.Range("H13").Value = 19999
$ws.Range("I13").Value = 9999
$ws.Range("K13").Value = 9999
$ws.Range("M13").Value = -9855
$ws.Range("H74").Value = 40002944
$ws.Range("I74").Value = 43481384
$ws.Range("K74").Value = 43481384
$ws.Range("M74").Value = -43480510
$ws.Range("H77").Value = 40002944
$ws.Range("I77").Value = 43481384
$ws.Range("K77").Value = 217406920
$ws.Range("M77").Value = -217402552
$ws.Range("H110").Value = 91840.27
$ws.Range("I110").Value = 100989.3
$ws.Range("J110").Value = 350
$ws.Range("K110").Value = 100989.3
$ws.Range("L110").Value = 350
$ws.Range("M110").Value = -98944.3
$ws.Range("N110").Value = -4440
$ws.Range("H116").Value = 3767.4
$ws.Range("I116").Value = 5594
$ws.Range("J116").Value = 2549.6667
$ws.Range("K116").Value = 5594
$ws.Range("L116").Value = 2549.6667
$ws.Range("M116").Value = -3300
$ws.Range("N116").Value = -7137.6667
$ws.Range("H122").Value = 4325.2666
$ws.Range("I122").Value = 2732.9
$ws.Range("K122").Value = 8198.700000000001
$ws.Range("M122").Value = -5748.700000000001
$ws.Range("H132").Value = 3126379.2
$ws.Range("I132").Value = 3126379.2
$ws.Range("K132").Value = 9379137.600000001
$ws.Range("M132").Value = -9376607.600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3767.4
$ws.Range("I3").Value = 5594
$ws.Range("J3").Value = 2549.6667
$ws.Range("K3").Value = 5594
$ws.Range("L3").Value = 2549.6667
$ws.Range("M3").Value = -5480
$ws.Range("N3").Value = -2777.6667
$ws.Range("H4").Value = 3244.25
$ws.Range("I4").Value = 1406
$ws.Range("K4").Value = 1406
$ws.Range("M4").Value = -1291

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 8684.75
$ws.Range("J7").Value = 299.66666
$ws.Range("L7").Value = 299.66666
$ws.Range("N7").Value = -525.66666
$ws.Range("H16").Value = 1554209
$ws.Range("I16").Value = 2718615.8
$ws.Range("K16").Value = 2718615.8
$ws.Range("M16").Value = -2718328.8
$ws.Range("H105").Value = 3402766.5
$ws.Range("I105").Value = 6803329.5
$ws.Range("J105").Value = 2203.6667
$ws.Range("K105").Value = 6803329.5
$ws.Range("L105").Value = 2203.6667
$ws.Range("M105").Value = -6801582.5
$ws.Range("N105").Value = -5697.6667
$ws.Range("H113").Value = 1554209
$ws.Range("I113").Value = 2718615.8
$ws.Range("K113").Value = 2718615.8
$ws.Range("M113").Value = -2716445.8
$ws.Range("H122").Value = 2626.72
$ws.Range("J122").Value = 2108.4443
$ws.Range("L122").Value = 6325.3329
$ws.Range("N122").Value = -11225.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4763922.5
$ws.Range("I137").Value = 11112887
$ws.Range("J137").Value = 2199.4167
$ws.Range("K137").Value = 33338661
$ws.Range("L137").Value = 6598.250100000001
$ws.Range("M137").Value = -33333561
$ws.Range("N137").Value = -16798.2501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H87").Value = 69999
$ws.Range("J87").Value = 69999
$ws.Range("L87").Value = 69999
$ws.Range("N87").Value = -72495
$ws.Range("H90").Value = 69999
$ws.Range("J90").Value = 69999
$ws.Range("L90").Value = 209997
$ws.Range("N90").Value = -222477
$ws.Range("H97").Value = 1111.3478
$ws.Range("I97").Value = 917.2222
$ws.Range("J97").Value = 1810.2
$ws.Range("K97").Value = 917.2222
$ws.Range("L97").Value = 1810.2
$ws.Range("M97").Value = -421.2222
$ws.Range("N97").Value = -2802.2
$ws.Range("H99").Value = 20924.75
$ws.Range("I99").Value = 4566.6665
$ws.Range("J99").Value = 69999
$ws.Range("K99").Value = 4566.6665
$ws.Range("L99").Value = 69999
$ws.Range("M99").Value = -2320.6665
$ws.Range("N99").Value = -74491
$ws.Range("H102").Value = 4072.125
$ws.Range("I102").Value = 3868.1428
$ws.Range("K102").Value = 3868.1428
$ws.Range("M102").Value = -2246.1428
$ws.Range("H126").Value = 2683.5
$ws.Range("J126").Value = 1800
$ws.Range("L126").Value = 5400
$ws.Range("N126").Value = -10340
$ws.Range("H132").Value = 8930757
$ws.Range("I132").Value = 9617566
$ws.Range("K132").Value = 28852698
$ws.Range("M132").Value = -28850168

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4995
$ws.Range("I7").Value = 4995
$ws.Range("K7").Value = 4995
$ws.Range("M7").Value = -4883
$ws.Range("H40").Value = 2934.5
$ws.Range("I40").Value = 3068.1428
$ws.Range("K40").Value = 3068.1428
$ws.Range("M40").Value = -2932.1428
$ws.Range("H126").Value = 4995
$ws.Range("I126").Value = 4995
$ws.Range("K126").Value = 14985
$ws.Range("M126").Value = -12515
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 219.71428
$ws.Range("J100").Value = 345
$ws.Range("L100").Value = 690
$ws.Range("N100").Value = -1772
$ws.Range("H122").Value = 2043.5
$ws.Range("I122").Value = 1933.2222
$ws.Range("J122").Value = 2185.2856
$ws.Range("K122").Value = 5799.6666
$ws.Range("L122").Value = 6555.8568
$ws.Range("M122").Value = -3349.6666
$ws.Range("N122").Value = -11455.8568
$ws.Range("H132").Value = 41677880
$ws.Range("I132").Value = 55561670
$ws.Range("J132").Value = 26496.666
$ws.Range("K132").Value = 166685010
$ws.Range("L132").Value = 79489.99800000001
$ws.Range("M132").Value = -166682480
$ws.Range("N132").Value = -84549.99800000001
